$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the text of the row that held the mangled, quote-embedded string,
#    replacing it with the clean version that survives into the final sheet.
$ws.Cells.Item(11, 1).Value = "!Buenos días! ¿Cómo va tu mañana?"

# 2. Remove the rows that were dropped entirely (delete bottom-up so the
#    remaining row numbers stay valid while we work).
$ws.Rows(23).EntireRow.Delete()   # "Hasta la próxima! No olvides evaluar..."
$ws.Rows(22).EntireRow.Delete()   # "No olvides evaluar el servicio..."
$ws.Rows(18).EntireRow.Delete()   # "Ya te echaba de menos..."
$ws.Rows(7).EntireRow.Delete()    # "Adelante, ¿en qué te puedo ayudar?"

# 3. Sort the remaining data (excluding the header row) ascending by CLASE
#    (column B): despedida, nombre, saludo.
$dataRange = $ws.Range("A2:B23")
$sortKey = $ws.Range("B2:B23")
$dataRange.Sort($sortKey, 1)

# 4. Re-apply the autofilter / sortState over the new, smaller range.
$ws.Range("A1:B23").AutoFilter()

Write-Host "done"
